# LV_TMTT0047023_VerifyFunctionalityOfRoundTripSectionOnCompanyInfoPage
# Add a new "FlagReason" worksheet (after "Warning") with Reason/Comment
# header row and a sample "Request to Change Company Type" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Add the new worksheet as the last tab (after "Warning")
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "FlagReason"

# Column widths (A ~32.44 chars, B = 56 chars)
$newSheet.Columns.Item(1).ColumnWidth = 31.67
$newSheet.Columns.Item(2).ColumnWidth = 56

# ---------------------------------------------------------------
# Header row (bold, centered horizontally + vertically)
# ---------------------------------------------------------------
$a1 = $newSheet.Range("A1")
$a1.Value = "Reason"
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4108

$b1 = $newSheet.Range("B1")
$b1.Value = "Comment"

# Re-use A1's exact direct formatting for B1 so both header cells share
# one style (avoids creating a stray intermediate style entry).
$a1.Copy() | Out-Null
$b1.PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------
# Data row
# ---------------------------------------------------------------
$a2 = $newSheet.Range("A2")
$a2.Value = "Request to Change Company Type"
$a2.HorizontalAlignment = -4131
$a2.VerticalAlignment = -4108

$b2 = $newSheet.Range("B2")
$b2.Value = "Requesting to either (i) change Company Type to Operating Company and Ownership to Private Equity Group or (ii) review the appropriateness of the round trip designation with CF operations"
$b2.VerticalAlignment = -4108
$b2.WrapText = $true

$newSheet.Rows.Item(2).RowHeight = 43.2

# Match the saved selection/active cell on the new sheet
$newSheet.Range("B12").Select() | Out-Null

$excel.CutCopyMode = $false
